# Append: 2025-11-16 01:55 JST
# Update the "取得日時" (acquired timestamp) column (A) for all existing
# data rows on the "ランサーズ" sheet from the previous run timestamp to
# the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-11-16 01:25:11"
$newTimestamp = "2025-11-16 01:55:14"

# Data rows start at row 2 (row 1 is the header) through row 16.
# (Note: compare using .Value2 -- .Value triggered an unrelated quirk in
# this host where the comparison always evaluated false.)
for ($r = 2; $r -le 16; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
